$d = $word.ActiveDocument

$replacements = @(
    @{old = "13×77=1001"; new = "40×61=2440"},
    @{old = "76×32=2432"; new = "64×32=2048"},
    @{old = "76×83=6308"; new = "76×87=6612"},
    @{old = "32×20=640";  new = "94×44=4136"},
    @{old = "71×92=6532"; new = "21×25=525"},
    @{old = "83×27=2241"; new = "62×26=1612"},
    @{old = "28×40=1120"; new = "92×47=4324"},
    @{old = "87×48=4176"; new = "79×45=3555"},
    @{old = "28×81=2268"; new = "52×66=3432"},
    @{old = "29×16=464";  new = "77×49=3773"},
    @{old = "92×63=5796"; new = "96×14=1344"},
    @{old = "66×75=4950"; new = "31×80=2480"},
    @{old = "84×27=2268"; new = "99×76=7524"},
    @{old = "53×32=1696"; new = "23×79=1817"},
    @{old = "18×99=1782"; new = "76×76=5776"},
    @{old = "38×86=3268"; new = "34×12=408"},
    @{old = "61×78=4758"; new = "15×49=735"},
    @{old = "27×42=1134"; new = "93×71=6603"},
    @{old = "36×85=3060"; new = "77×81=6237"},
    @{old = "24×61=1464"; new = "90×65=5850"},
    @{old = "69×43=2967"; new = "30×51=1530"},
    @{old = "43×99=4257"; new = "15×97=1455"},
    @{old = "28×45=1260"; new = "47×40=1880"},
    @{old = "39×67=2613"; new = "92×18=1656"},
    @{old = "77×39=3003"; new = "80×35=2800"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
